$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches source formatting).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.248.20'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.263.02'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '495.58'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.62'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  -0.93%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.83'
$ws.Range('E12').Value = '  +4.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.92'
$ws.Range('E13').Value = '  +5.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.664.49'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.227.93'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.269.16'
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.22'
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '302.97'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '60.70'
$ws.Range('E23').Value = '  -3.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.996'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('E26').Value = '  +3.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.52'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.97'
$ws.Range('E29').Value = '  +2.04%  '
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.939'
$ws.Range('E35').Value = '  +6.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.20'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.80'
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '124.71'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0893'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '241.39'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.10'
$ws.Range('E50').Value = '  -1.12%  '
$ws.Range('E51').Value = '  -0.42%  '
